$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 850
$ws.Range("J10").Value = 850
$ws.Range("L10").Value = 850
$ws.Range("N10").Value = -1436
$ws.Range("H40").Value = 5266
$ws.Range("I40").Value = 4800
$ws.Range("J40").Value = 5499
$ws.Range("K40").Value = 4800
$ws.Range("L40").Value = 5499
$ws.Range("M40").Value = -4625
$ws.Range("N40").Value = -5849
$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 5000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 25000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -31240
$ws.Range("H76").Value = 4100.3335
$ws.Range("I76").Value = 4000.3333
$ws.Range("K76").Value = 4000.3333
$ws.Range("M76").Value = -3685.3333
$ws.Range("H79").Value = 4100.3335
$ws.Range("I79").Value = 4000.3333
$ws.Range("K79").Value = 4000.3333
$ws.Range("M79").Value = -2908.3333
$ws.Range("H92").Value = 174.125
$ws.Range("J92").Value = 271.25
$ws.Range("L92").Value = 271.25
$ws.Range("N92").Value = -2767.25
$ws.Range("H96").Value = 2639.9167
$ws.Range("I96").Value = 690
$ws.Range("J96").Value = 2817.182
$ws.Range("K96").Value = 2070
$ws.Range("L96").Value = 8451.545999999998
$ws.Range("M96").Value = -697
$ws.Range("N96").Value = -11197.546
$ws.Range("H100").Value = 1573.9524
$ws.Range("I100").Value = 1521.7142
$ws.Range("J100").Value = 1678.4286
$ws.Range("K100").Value = 1521.7142
$ws.Range("L100").Value = 1678.4286
$ws.Range("M100").Value = -980.7141999999999
$ws.Range("N100").Value = -2760.4286
$ws.Range("H111").Value = 397
$ws.Range("I111").Value = 352
$ws.Range("K111").Value = 1056
$ws.Range("M111").Value = 2011
$ws.Range("H129").Value = 1177603
$ws.Range("I129").Value = 770156.0600000001
$ws.Range("K129").Value = 2310468.18
$ws.Range("M129").Value = -2305468.18
$ws.Range("H137").Value = 10535382
$ws.Range("I137").Value = 22225198
$ws.Range("J137").Value = 14547.5
$ws.Range("K137").Value = 66675594
$ws.Range("L137").Value = 43642.5
$ws.Range("M137").Value = -66673044
$ws.Range("N137").Value = -48742.5
$ws.Range("H138").Value = 5951.92
$ws.Range("I138").Value = 3681.1177
$ws.Range("J138").Value = 7121.727
$ws.Range("K138").Value = 11043.3531
$ws.Range("L138").Value = 21365.181
$ws.Range("M138").Value = -5903.3531
$ws.Range("N138").Value = -31645.181
$ws.Range("H141").Value = 1687.174
$ws.Range("I141").Value = 1622.9546
$ws.Range("J141").Value = 3100
$ws.Range("K141").Value = 4868.8638
$ws.Range("L141").Value = 9300
$ws.Range("M141").Value = 311.1361999999999
$ws.Range("N141").Value = -19660

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 4517
$ws.Range("J25").Value = 1018
$ws.Range("L25").Value = 1018
$ws.Range("N25").Value = -1822
$ws.Range("H32").Value = 1207412
$ws.Range("I32").Value = 1540409.9
$ws.Range("K32").Value = 1540409.9
$ws.Range("M32").Value = -1540122.9
$ws.Range("H35").Value = 1900
$ws.Range("I35").Value = 1900
$ws.Range("K35").Value = 1900
$ws.Range("M35").Value = -1494
$ws.Range("H41").Value = 2550.5557
$ws.Range("I41").Value = 1422.1428
$ws.Range("K41").Value = 1422.1428
$ws.Range("M41").Value = -1008.1428
$ws.Range("H45").Value = 3061.68
$ws.Range("I45").Value = 2733.7896
$ws.Range("K45").Value = 2733.7896
$ws.Range("M45").Value = -2356.7896
$ws.Range("H57").Value = 8996
$ws.Range("I57").Value = 8996
$ws.Range("K57").Value = 8996
$ws.Range("M57").Value = -8512
$ws.Range("H61").Value = 56581816
$ws.Range("I61").Value = 45954544
$ws.Range("K61").Value = 45954544
$ws.Range("M61").Value = -45954332
$ws.Range("H74").Value = 1285263.9
$ws.Range("I74").Value = 1815152.1
$ws.Range("K74").Value = 1815152.1
$ws.Range("M74").Value = -1814278.1
$ws.Range("H77").Value = 1285263.9
$ws.Range("I77").Value = 1815152.1
$ws.Range("K77").Value = 9075760.5
$ws.Range("M77").Value = -9071392.5
$ws.Range("H102").Value = 1268.9166
$ws.Range("I102").Value = 1170.9688
$ws.Range("K102").Value = 1170.9688
$ws.Range("M102").Value = 451.0311999999999
$ws.Range("H132").Value = 5408.4717
$ws.Range("I132").Value = 3243.875
$ws.Range("K132").Value = 9731.625
$ws.Range("M132").Value = -7201.625
$ws.Range("H135").Value = 101447.75
$ws.Range("J135").Value = 101447.75
$ws.Range("L135").Value = 101447.75
$ws.Range("N135").Value = -111587.75
$ws.Range("H136").Value = 56581816
$ws.Range("I136").Value = 45954544
$ws.Range("K136").Value = 137863632
$ws.Range("M136").Value = -137861082

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H17").Value = 6200
$ws.Range("J17").Value = 6200
$ws.Range("L17").Value = 6200
$ws.Range("N17").Value = -6544
$ws.Range("H107").Value = 1258.7142
$ws.Range("I107").Value = 1258.7142
$ws.Range("K107").Value = 1258.7142
$ws.Range("M107").Value = 661.2858000000001
$ws.Range("H134").Value = 15115185
$ws.Range("I134").Value = 13706841
$ws.Range("K134").Value = 41120523
$ws.Range("M134").Value = -41117988
$ws.Range("H135").Value = 65921.44500000001
$ws.Range("J135").Value = 65921.44500000001
$ws.Range("L135").Value = 65921.44500000001
$ws.Range("N135").Value = -76061.44500000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1460.8182
$ws.Range("I22").Value = 854.5
$ws.Range("K22").Value = 854.5
$ws.Range("M22").Value = -504.5
$ws.Range("H31").Value = 421356.5
$ws.Range("I31").Value = 909691.9
$ws.Range("J31").Value = 4835.147
$ws.Range("K31").Value = 909691.9
$ws.Range("L31").Value = 4835.147
$ws.Range("M31").Value = -909396.9
$ws.Range("N31").Value = -5425.147
$ws.Range("H34").Value = 421356.5
$ws.Range("I34").Value = 909691.9
$ws.Range("J34").Value = 4835.147
$ws.Range("K34").Value = 909691.9
$ws.Range("L34").Value = 4835.147
$ws.Range("M34").Value = -909489.9
$ws.Range("N34").Value = -5239.147
$ws.Range("H58").Value = 23354382
$ws.Range("I58").Value = 41673416
$ws.Range("K58").Value = 41673416
$ws.Range("M58").Value = -41673213
$ws.Range("H62").Value = 4399
$ws.Range("I62").Value = 4638.6
$ws.Range("J62").Value = 3999.6667
$ws.Range("K62").Value = 4638.6
$ws.Range("L62").Value = 3999.6667
$ws.Range("M62").Value = -4014.6
$ws.Range("N62").Value = -5247.6667
$ws.Range("H65").Value = 4399
$ws.Range("I65").Value = 4638.6
$ws.Range("J65").Value = 3999.6667
$ws.Range("K65").Value = 23193
$ws.Range("L65").Value = 19998.3335
$ws.Range("M65").Value = -20073
$ws.Range("N65").Value = -26238.3335
$ws.Range("H94").Value = 21298.4
$ws.Range("J94").Value = 1623
$ws.Range("L94").Value = 1623
$ws.Range("N94").Value = -2525
$ws.Range("H96").Value = 30666.334
$ws.Range("J96").Value = 30666.334
$ws.Range("L96").Value = 30666.334
$ws.Range("N96").Value = -36158.334
$ws.Range("H122").Value = 13838.823
$ws.Range("I122").Value = 1954.2858
$ws.Range("K122").Value = 5862.857400000001
$ws.Range("M122").Value = -3412.857400000001
$ws.Range("H132").Value = 5617.706
$ws.Range("I132").Value = 5549.846
$ws.Range("K132").Value = 16649.538
$ws.Range("M132").Value = -14119.538
$ws.Range("H134").Value = 4285.886
$ws.Range("I134").Value = 2059.9443
$ws.Range("K134").Value = 6179.8329
$ws.Range("M134").Value = -3644.8329
$ws.Range("H136").Value = 23354382
$ws.Range("I136").Value = 41673416
$ws.Range("K136").Value = 125020248
$ws.Range("M136").Value = -125017698

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 100.13333
$ws.Range("I2").Value = 39.857143
$ws.Range("J2").Value = 152.875
$ws.Range("K2").Value = 239.142858
$ws.Range("L2").Value = 917.25
$ws.Range("M2").Value = -126.142858
$ws.Range("N2").Value = -1143.25
$ws.Range("H4").Value = 3737035.2
$ws.Range("I4").Value = 4057335.2
$ws.Range("J4").Value = 199
$ws.Range("K4").Value = 12172005.6
$ws.Range("L4").Value = 597
$ws.Range("M4").Value = -12171893.6
$ws.Range("N4").Value = -821
$ws.Range("H110").Value = 2998.5
$ws.Range("I110").Value = 2998.5
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 8995.5
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -4905.5
$ws.Range("N110").ClearContents()
$ws.Range("H117").Value = 1208.375
$ws.Range("I117").Value = 272.5
$ws.Range("J117").Value = 4016
$ws.Range("K117").Value = 817.5
$ws.Range("L117").Value = 12048
$ws.Range("M117").Value = 2624.5
$ws.Range("N117").Value = -18932
$ws.Range("H120").Value = 8500
$ws.Range("I120").Value = 8500
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 25500
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -20662
$ws.Range("N120").ClearContents()
$ws.Range("H132").Value = 1108.1
$ws.Range("I132").Value = 556.4
$ws.Range("J132").Value = 1659.8
$ws.Range("K132").Value = 5007.599999999999
$ws.Range("L132").Value = 14938.2
$ws.Range("M132").Value = -2477.599999999999
$ws.Range("N132").Value = -19998.2

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 8000
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H102").Value = 1410.4193
$ws.Range("I102").Value = 1300.3846
$ws.Range("K102").Value = 1300.3846
$ws.Range("M102").Value = 321.6153999999999
$ws.Range("H107").Value = 800.1667
$ws.Range("I107").Value = 840.3
$ws.Range("K107").Value = 840.3
$ws.Range("M107").Value = 1079.7
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H126").Value = 3085.0908
$ws.Range("J126").Value = 3833.4
$ws.Range("L126").Value = 11500.2
$ws.Range("N126").Value = -16440.2
$ws.Range("H132").Value = 20068.863
$ws.Range("I132").Value = 27585.715
$ws.Range("K132").Value = 82757.145
$ws.Range("M132").Value = -80227.145

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5497.8335
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H46").Value = 3250.3333
$ws.Range("I46").Value = 2794.6
$ws.Range("K46").Value = 2794.6
$ws.Range("M46").Value = -2606.6
$ws.Range("H68").Value = 5500.5
$ws.Range("I68").Value = 5500.5
$ws.Range("K68").Value = 5500.5
$ws.Range("M68").Value = -4751.5
$ws.Range("H71").Value = 5500.5
$ws.Range("I71").Value = 5500.5
$ws.Range("K71").Value = 27502.5
$ws.Range("M71").Value = -23758.5
$ws.Range("H107").Value = 5500
$ws.Range("I107").Value = 5500
$ws.Range("K107").Value = 5500
$ws.Range("M107").Value = -3580
$ws.Range("H122").Value = 5907.4375
$ws.Range("I122").Value = 6552.727
$ws.Range("K122").Value = 19658.181
$ws.Range("M122").Value = -17208.181
$ws.Range("H126").Value = 5497.8335
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 2225811
$ws.Range("I132").Value = 3336871.5
$ws.Range("K132").Value = 10010614.5
$ws.Range("M132").Value = -10008084.5
$ws.Range("H136").Value = 19178558
$ws.Range("I136").Value = 13902101
$ws.Range("K136").Value = 41706303
$ws.Range("M136").Value = -41703753

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 19621
$ws.Range("I62").Value = 11741.25
$ws.Range("J62").Value = 27500.75
$ws.Range("K62").Value = 11741.25
$ws.Range("L62").Value = 27500.75
$ws.Range("M62").Value = -11117.25
$ws.Range("N62").Value = -28748.75
$ws.Range("H65").Value = 19621
$ws.Range("I65").Value = 11741.25
$ws.Range("J65").Value = 27500.75
$ws.Range("K65").Value = 58706.25
$ws.Range("L65").Value = 137503.75
$ws.Range("M65").Value = -55586.25
$ws.Range("N65").Value = -143743.75
$ws.Range("H113").Value = 1372.5714
$ws.Range("I113").Value = 718.8333
$ws.Range("K113").Value = 2156.4999
$ws.Range("M113").Value = 13.5001000000002
$ws.Range("H126").Value = 3635.8096
$ws.Range("I126").Value = 2413.353
$ws.Range("K126").Value = 7240.059
$ws.Range("M126").Value = -4770.059
$ws.Range("H132").Value = 5749610.5
$ws.Range("I132").Value = 7248818
$ws.Range("J132").Value = 2647.1667
$ws.Range("K132").Value = 21746454
$ws.Range("L132").Value = 7941.500100000001
$ws.Range("M132").Value = -21743924
$ws.Range("N132").Value = -13001.5001
